$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (pushes existing rows 2-23 down to 3-24)
$ws.Rows.Item(2).Insert()

# The inserted row inherits formatting from the row above; clear it so the
# new data row matches the unstyled look of the other data rows.
$ws.Rows.Item(2).ClearFormats()

# New top row: today's date with same price values as the rest of the series.
# Force the date cell to be stored as literal text (not auto-parsed into a
# date serial number), then drop back to the Normal style so no stray
# number-format style lingers on the cell.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2025-12-13"
$ws.Cells.Item(2, 1).Style = "Normal"
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
